# Rename the "firework" sheet to "yew" and remove the second data row
# (the row that held the "Pursuit of Happiness ... [Extended Steve Aoki Remix]"
# entry), leaving only the header row behind.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("firework")
$ws.Name = "yew"

# Delete the whole second row (shifts nothing else up since it's the last row)
$ws.Rows.Item(2).Delete()
